$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 46, shifting the existing
# rows 46:55 down to 47:56 (all of their values travel with them).
$ws.Rows("46:46").Insert()

# Populate the newly-inserted row 46 with the new weekly record.
$ws.Cells.Item(46, 1).Value = 7
$ws.Cells.Item(46, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(46, 3).Value = "Ñuble"
$ws.Cells.Item(46, 4).Value = 44873
$ws.Cells.Item(46, 5).Value = 16
$ws.Cells.Item(46, 6).Value = 100112001
$ws.Cells.Item(46, 7).Value = "Berenjena"
$ws.Cells.Item(46, 8).Value = "Sin especificar"
$ws.Cells.Item(46, 9).Value = "Primera"
$ws.Cells.Item(46, 10).Value = 60
$ws.Cells.Item(46, 11).Value = 13000
$ws.Cells.Item(46, 12).Value = 14000
$ws.Cells.Item(46, 13).Value = 13500
$ws.Cells.Item(46, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(46, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(46, 16).Value = 225
$ws.Cells.Item(46, 17).Value = 60
$ws.Cells.Item(46, 18).Value = "Hortaliza"
